$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 143, shifting existing rows 143-188 down to 146-191.
$ws.Range("A143:A145").EntireRow.Insert()

# New row 143: Chirimoya - Especial - Provincia de Limarí, new sampling date 44845
$ws.Cells.Item(143,1).Value = 8
$ws.Cells.Item(143,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143,3).Value = "Coquimbo"
$ws.Cells.Item(143,4).Value = 44845
$ws.Cells.Item(143,5).Value = 4
$ws.Cells.Item(143,6).Value = "Fruta"
$ws.Cells.Item(143,7).Value = 100107
$ws.Cells.Item(143,8).Value = "Otros"
$ws.Cells.Item(143,9).Value = 100107002
$ws.Cells.Item(143,10).Value = "Chirimoya"
$ws.Cells.Item(143,11).Value = "Cultivar IV Región"
$ws.Cells.Item(143,12).Value = "Especial"
$ws.Cells.Item(143,13).Value = 400
$ws.Cells.Item(143,14).Value = 21000
$ws.Cells.Item(143,15).Value = 22000
$ws.Cells.Item(143,16).Value = 21500
$ws.Cells.Item(143,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(143,18).Value = "Provincia de Limarí"
$ws.Cells.Item(143,19).Value = 2150
$ws.Cells.Item(143,20).Value = 10

# New row 144: Chirimoya - Primera - Provincia de Limarí, new sampling date 44845
$ws.Cells.Item(144,1).Value = 8
$ws.Cells.Item(144,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(144,3).Value = "Coquimbo"
$ws.Cells.Item(144,4).Value = 44845
$ws.Cells.Item(144,5).Value = 4
$ws.Cells.Item(144,6).Value = "Fruta"
$ws.Cells.Item(144,7).Value = 100107
$ws.Cells.Item(144,8).Value = "Otros"
$ws.Cells.Item(144,9).Value = 100107002
$ws.Cells.Item(144,10).Value = "Chirimoya"
$ws.Cells.Item(144,11).Value = "Cultivar IV Región"
$ws.Cells.Item(144,12).Value = "Primera"
$ws.Cells.Item(144,13).Value = 300
$ws.Cells.Item(144,14).Value = 18000
$ws.Cells.Item(144,15).Value = 19000
$ws.Cells.Item(144,16).Value = 18500
$ws.Cells.Item(144,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(144,18).Value = "Provincia de Limarí"
$ws.Cells.Item(144,19).Value = 1850
$ws.Cells.Item(144,20).Value = 10

# New row 145: Chirimoya - Segunda - Provincia de Limarí, new sampling date 44845
$ws.Cells.Item(145,1).Value = 8
$ws.Cells.Item(145,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(145,3).Value = "Coquimbo"
$ws.Cells.Item(145,4).Value = 44845
$ws.Cells.Item(145,5).Value = 4
$ws.Cells.Item(145,6).Value = "Fruta"
$ws.Cells.Item(145,7).Value = 100107
$ws.Cells.Item(145,8).Value = "Otros"
$ws.Cells.Item(145,9).Value = 100107002
$ws.Cells.Item(145,10).Value = "Chirimoya"
$ws.Cells.Item(145,11).Value = "Cultivar IV Región"
$ws.Cells.Item(145,12).Value = "Segunda"
$ws.Cells.Item(145,13).Value = 240
$ws.Cells.Item(145,14).Value = 15000
$ws.Cells.Item(145,15).Value = 16000
$ws.Cells.Item(145,16).Value = 15500
$ws.Cells.Item(145,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(145,18).Value = "Provincia de Limarí"
$ws.Cells.Item(145,19).Value = 1550
$ws.Cells.Item(145,20).Value = 10
